$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 558
$ws1.Range("F13").Value = 166
$ws1.Range("F18").Value = 5104
$ws1.Range("F19").Value = 57
$ws1.Range("F20").Value = 837
$ws1.Range("F22").Value = 2276

# Sheet "全部类型" (fourth sheet) - update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 558
$ws4.Range("F13").Value = 166
$ws4.Range("F18").Value = 5104
$ws4.Range("F20").Value = 57
$ws4.Range("F22").Value = 837
$ws4.Range("F24").Value = 2276
